$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10464.889
$ws.Range("J51").Value = 10464.889
$ws.Range("L51").Value = 10464.889
$ws.Range("N51").Value = -11432.889
$ws.Range("H93").Value = 30941.709
$ws.Range("J93").Value = 30941.709
$ws.Range("L93").Value = 30941.709
$ws.Range("N93").Value = -35933.709
$ws.Range("H115").Value = 1141.9231
$ws.Range("I115").Value = 1141.9231
$ws.Range("K115").Value = 3425.7693
$ws.Range("M115").Value = -1858.7693
$ws.Range("H129").Value = 856.55
$ws.Range("J129").Value = 871.3917
$ws.Range("L129").Value = 2614.1751
$ws.Range("N129").Value = -12614.1751
$ws.Range("H131").Value = 3523.5557
$ws.Range("I131").Value = 2804.875
$ws.Range("K131").Value = 8414.625
$ws.Range("M131").Value = -3374.625
$ws.Range("H137").Value = 1907004.6
$ws.Range("I137").Value = 2646900.8
$ws.Range("J137").Value = 4414.7144
$ws.Range("K137").Value = 7940702.399999999
$ws.Range("L137").Value = 13244.1432
$ws.Range("M137").Value = -7938152.399999999
$ws.Range("N137").Value = -18344.1432
$ws.Range("H138").Value = 2715.5
$ws.Range("J138").Value = 2865.8462
$ws.Range("L138").Value = 8597.5386
$ws.Range("N138").Value = -18877.5386
$ws.Range("H141").Value = 84917.086
$ws.Range("I141").Value = 92191.37
$ws.Range("K141").Value = 276574.11
$ws.Range("M141").Value = -271394.11

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 439.85715
$ws.Range("J2").Value = 649
$ws.Range("L2").Value = 649
$ws.Range("H32").Value = 6066.2095
$ws.Range("I32").Value = 4895.6055
$ws.Range("J32").Value = 14962.8
$ws.Range("K32").Value = 4895.6055
$ws.Range("L32").Value = 14962.8
$ws.Range("M32").Value = -4608.6055
$ws.Range("N32").Value = -15536.8
$ws.Range("H61").Value = 2179.2
$ws.Range("I61").Value = 1519.6
$ws.Range("J61").Value = 2838.8
$ws.Range("K61").Value = 1519.6
$ws.Range("L61").Value = 2838.8
$ws.Range("M61").Value = -1307.6
$ws.Range("N61").Value = -3262.8
$ws.Range("H116").Value = 439.85715
$ws.Range("J116").Value = 649
$ws.Range("L116").Value = 649
$ws.Range("H132").Value = 1609.4736
$ws.Range("I132").Value = 800.34485
$ws.Range("J132").Value = 4216.6665
$ws.Range("K132").Value = 2401.03455
$ws.Range("L132").Value = 12649.9995
$ws.Range("M132").Value = 128.9654500000001
$ws.Range("N132").Value = -17709.9995
$ws.Range("H135").Value = 50095.668
$ws.Range("J135").Value = 50095.668
$ws.Range("L135").Value = 50095.668
$ws.Range("H136").Value = 2179.2
$ws.Range("I136").Value = 1519.6
$ws.Range("J136").Value = 2838.8
$ws.Range("K136").Value = 4558.799999999999
$ws.Range("L136").Value = 8516.400000000001
$ws.Range("M136").Value = -2008.799999999999
$ws.Range("N136").Value = -13616.4
$ws.Range("H137").Value = 40872.5
$ws.Range("J137").Value = 40872.5
$ws.Range("L137").Value = 40872.5
$ws.Range("N137").Value = -51072.5
$ws.Range("N2").Value = -875
$ws.Range("N116").Value = -5237
$ws.Range("N135").Value = -60235.668

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 439.85715
$ws.Range("J3").Value = 649
$ws.Range("L3").Value = 649
$ws.Range("H95").Value = 33933.332
$ws.Range("J95").Value = 33933.332
$ws.Range("L95").Value = 33933.332
$ws.Range("N95").Value = -39425.332
$ws.Range("H134").Value = 2430.121
$ws.Range("J134").Value = 4091.25
$ws.Range("L134").Value = 12273.75
$ws.Range("N134").Value = -17343.75
$ws.Range("H137").Value = 45720
$ws.Range("J137").Value = 45720
$ws.Range("L137").Value = 45720
$ws.Range("N137").Value = -55920
$ws.Range("N3").Value = -877

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3362.7368
$ws.Range("I58").Value = 1935.8182
$ws.Range("J58").Value = 8192.308000000001
$ws.Range("K58").Value = 1935.8182
$ws.Range("L58").Value = 8192.308000000001
$ws.Range("M58").Value = -1732.8182
$ws.Range("N58").Value = -8598.308000000001
$ws.Range("H86").Value = 2306
$ws.Range("I86").Value = 2531
$ws.Range("J86").Value = 2171
$ws.Range("K86").Value = 2531
$ws.Range("L86").Value = 2171
$ws.Range("M86").Value = -1408
$ws.Range("N86").Value = -4417
$ws.Range("H89").Value = 2306
$ws.Range("I89").Value = 2531
$ws.Range("J89").Value = 2171
$ws.Range("K89").Value = 12655
$ws.Range("L89").Value = 10855
$ws.Range("M89").Value = -7039
$ws.Range("N89").Value = -22087
$ws.Range("H122").Value = 15000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("H136").Value = 3362.7368
$ws.Range("I136").Value = 1935.8182
$ws.Range("J136").Value = 8192.308000000001
$ws.Range("K136").Value = 5807.4546
$ws.Range("L136").Value = 24576.924
$ws.Range("M136").Value = -3257.4546
$ws.Range("N136").Value = -29676.924
$ws.Range("M122").ClearContents()
$ws.Range("N123").ClearContents()
$ws.Range("N125").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3826.5
$ws.Range("I3").Value = 3349.75
$ws.Range("K3").Value = 10049.25
$ws.Range("M3").Value = -9937.25
$ws.Range("H38").Value = 138.53334
$ws.Range("I38").Value = 99.833336
$ws.Range("J38").Value = 164.33333
$ws.Range("K38").Value = 299.500008
$ws.Range("L38").Value = 492.99999
$ws.Range("M38").Value = 47.49999200000002
$ws.Range("N38").Value = -1186.99999
$ws.Range("H113").Value = 2717981.5
$ws.Range("I113").Value = 601.8333
$ws.Range("J113").Value = 5682395.5
$ws.Range("K113").Value = 1805.4999
$ws.Range("L113").Value = 17047186.5
$ws.Range("M113").Value = 364.5001
$ws.Range("N113").Value = -17051526.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 19083.533
$ws.Range("I43").Value = 1100
$ws.Range("J43").Value = 28075.3
$ws.Range("K43").Value = 1100
$ws.Range("L43").Value = 28075.3
$ws.Range("M43").Value = -949
$ws.Range("N43").Value = -28377.3
$ws.Range("H102").Value = 2008.3125
$ws.Range("I102").Value = 1321
$ws.Range("J102").Value = 4986.6665
$ws.Range("K102").Value = 1321
$ws.Range("L102").Value = 4986.6665
$ws.Range("M102").Value = 301
$ws.Range("N102").Value = -8230.666499999999
$ws.Range("H113").Value = 1729.9412
$ws.Range("I113").Value = 1750.75
$ws.Range("J113").Value = 1680
$ws.Range("K113").Value = 1750.75
$ws.Range("L113").Value = 1680
$ws.Range("M113").Value = 419.25
$ws.Range("N113").Value = -6020
$ws.Range("H132").Value = 3528.4092
$ws.Range("I132").Value = 1559.3334
$ws.Range("J132").Value = 5891.3
$ws.Range("K132").Value = 4678.0002
$ws.Range("L132").Value = 17673.9
$ws.Range("M132").Value = -2148.0002
$ws.Range("N132").Value = -22733.9
$ws.Range("H137").Value = 38770
$ws.Range("J137").Value = 38770
$ws.Range("L137").Value = 38770
$ws.Range("N137").Value = -48970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 23000
$ws.Range("J3").Value = 23000
$ws.Range("L3").Value = 23000
$ws.Range("H15").Value = 23000
$ws.Range("J15").Value = 23000
$ws.Range("L15").Value = 23000
$ws.Range("H112").Value = 31710.525
$ws.Range("J112").Value = 31710.525
$ws.Range("L112").Value = 31710.525
$ws.Range("N112").Value = -34664.525
$ws.Range("H132").Value = 5761.6924
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5908.5
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 17725.5
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -22785.5
$ws.Range("H136").Value = 4220.16
$ws.Range("I136").Value = 1323.3846
$ws.Range("J136").Value = 7358.3335
$ws.Range("K136").Value = 3970.1538
$ws.Range("L136").Value = 22075.0005
$ws.Range("M136").Value = -1420.1538
$ws.Range("N136").Value = -27175.0005
$ws.Range("N3").Value = -23224
$ws.Range("N15").Value = -23340

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 28857.143
$ws.Range("J86").Value = 28857.143
$ws.Range("L86").Value = 28857.143
$ws.Range("N86").Value = -31103.143
$ws.Range("H89").Value = 28857.143
$ws.Range("J89").Value = 28857.143
$ws.Range("L89").Value = 144285.715
$ws.Range("N89").Value = -155517.715
